# Insert a new data row at row 226 (pushing existing rows 226..279 down to 227..280)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 226; this shifts row 226 onward down by one
# and the sheet dimension grows from A1:T279 to A1:T280 automatically.
$ws.Rows.Item(226).Insert()

# Populate the newly inserted (now empty) row 226 with the new record.
$ws.Cells.Item(226, 1).Value  = 11
$ws.Cells.Item(226, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(226, 3).Value  = "Bíobío"
$ws.Cells.Item(226, 4).Value  = 45211
$ws.Cells.Item(226, 5).Value  = 8
$ws.Cells.Item(226, 6).Value  = "Fruta"
$ws.Cells.Item(226, 7).Value  = 100102
$ws.Cells.Item(226, 8).Value  = "Cítricos"
$ws.Cells.Item(226, 9).Value  = 100102004
$ws.Cells.Item(226, 10).Value = "Mandarina"
$ws.Cells.Item(226, 11).Value = "Murcott"
$ws.Cells.Item(226, 12).Value = "Primera"
$ws.Cells.Item(226, 13).Value = 200
$ws.Cells.Item(226, 14).Value = 9000
$ws.Cells.Item(226, 15).Value = 10000
$ws.Cells.Item(226, 16).Value = 9500
$ws.Cells.Item(226, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(226, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(226, 19).Value = 528
$ws.Cells.Item(226, 20).Value = 18
